$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = "12 A"
$ws.Range("O3").Value = ""
$ws.Range("O4").Value = ""
$ws.Range("O5").Value = "12 A"
